$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$d = $s.Design
Write-Host "Current Name:" $d.Name
$d.Name = "Office Theme"
Write-Host "done"
